$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the Actual Outcome / Fail-Pass for test case #4 (row 5)
$ws.Range("F5").Value = "Same as expected outcome."
$ws.Range("G5").Value = "Pass"

# Update the active selection to match the saved view state
$ws.Range("G4:G5").Select()
